$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.477.53"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.289.57"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").Value = "2.700.22"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.26%  "
$ws.Range("D15").Value = "54.438.29"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "2.292.87"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "304.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("E26").Value = "  +3.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.967"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.88%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "126.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0898"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "242.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.374"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("E51").Value = "  -0.53%  "
